$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2175324675324675
$ws.Range("C2").Value = 0.5097402597402597
$ws.Range("J2").Value = 0.01948051948051948
$ws.Range("P2").Value = 0.1525974025974026
$ws.Range("S2").Value = 0.1006493506493507
$ws.Range("C3").Value = 0.01197604790419162
$ws.Range("J3").Value = 0.05389221556886228
$ws.Range("P3").Value = 0.7724550898203593
$ws.Range("S3").Value = 0.1616766467065868
$ws.Range("J4").Value = 0.119047619047619
$ws.Range("P4").Value = 0.7380952380952381
$ws.Range("S4").Value = 0.1428571428571428
$ws.Range("B6").Value = 0.04032258064516129
$ws.Range("D6").Value = 0.008064516129032258
$ws.Range("F6").Value = 0.06451612903225806
$ws.Range("J6").Value = 0.2862903225806452
$ws.Range("O6").Value = 0.004032258064516129
$ws.Range("Q6").Value = 0.1532258064516129
$ws.Range("R6").Value = 0.1008064516129032
$ws.Range("S6").Value = 0.3427419354838709
$ws.Range("B7").Value = 0.09134615384615384
$ws.Range("D7").Value = 0.01442307692307692
$ws.Range("F7").Value = 0.04326923076923077
$ws.Range("J7").Value = 0.1394230769230769
$ws.Range("O7").Value = 0.009615384615384616
$ws.Range("Q7").Value = 0.1298076923076923
$ws.Range("R7").Value = 0.125
$ws.Range("S7").Value = 0.4471153846153846
$ws.Range("B8").Value = 0.09513274336283185
$ws.Range("D8").Value = 0.01991150442477876
$ws.Range("F8").Value = 0.07964601769911504
$ws.Range("J8").Value = 0.1438053097345133
$ws.Range("O8").Value = 0.02433628318584071
$ws.Range("Q8").Value = 0.1393805309734513
$ws.Range("R8").Value = 0.1305309734513274
$ws.Range("S8").Value = 0.3672566371681416
$ws.Range("B9").Value = 0.154320987654321
$ws.Range("D9").Value = 0.01851851851851852
$ws.Range("F9").Value = 0.07407407407407407
$ws.Range("J9").Value = 0.1481481481481481
$ws.Range("O9").Value = 0.01851851851851852
$ws.Range("Q9").Value = 0.1296296296296296
$ws.Range("R9").Value = 0.1172839506172839
$ws.Range("S9").Value = 0.3395061728395062
$ws.Range("B10").Value = 0.09593604263824117
$ws.Range("D10").Value = 0.0173217854763491
$ws.Range("E10").Value = 0.0006662225183211193
$ws.Range("F10").Value = 0.06395736175882745
$ws.Range("J10").Value = 0.128580946035976
$ws.Range("O10").Value = 0.02664890073284477
$ws.Range("Q10").Value = 0.2338441039307129
$ws.Range("R10").Value = 0.09926715522984676
$ws.Range("S10").Value = 0.3337774816788808
$ws.Range("G11").Value = 0.10580204778157
$ws.Range("J11").Value = 0.1228668941979522
$ws.Range("K11").Value = 0.174061433447099
$ws.Range("L11").Value = 0.5938566552901023
$ws.Range("S11").Value = 0.003412969283276451
$ws.Range("G12").Value = 0.8146067415730337
$ws.Range("J12").Value = 0.1573033707865168
$ws.Range("L12").Value = 0.01123595505617977
$ws.Range("S12").Value = 0.01685393258426966
$ws.Range("F15").Value = 0.02066115702479339
$ws.Range("H15").Value = 0.1818181818181818
$ws.Range("I15").Value = 0.01239669421487603
$ws.Range("J15").Value = 0.3223140495867768
$ws.Range("K15").Value = 0.07024793388429752
$ws.Range("M15").Value = 0.01652892561983471
$ws.Range("N15").Value = 0.004132231404958678
$ws.Range("O15").Value = 0.08677685950413223
$ws.Range("S15").Value = 0.2851239669421488
$ws.Range("F16").Value = 0.015
$ws.Range("H16").Value = 0.22
$ws.Range("I16").Value = 0.03
$ws.Range("J16").Value = 0.47
$ws.Range("K16").Value = 0.09
$ws.Range("M16").Value = 0.025
$ws.Range("N16").Value = 0.005
$ws.Range("O16").Value = 0.04
$ws.Range("S16").Value = 0.105
$ws.Range("F17").Value = 0.0280561122244489
$ws.Range("H17").Value = 0.156312625250501
$ws.Range("I17").Value = 0.08216432865731463
$ws.Range("J17").Value = 0.4589178356713427
$ws.Range("K17").Value = 0.09218436873747494
$ws.Range("M17").Value = 0.02404809619238477
$ws.Range("N17").Value = 0.002004008016032064
$ws.Range("O17").Value = 0.05210420841683366
$ws.Range("S17").Value = 0.1042084168336673
$ws.Range("F18").Value = 0.01818181818181818
$ws.Range("H18").Value = 0.1818181818181818
$ws.Range("I18").Value = 0.08727272727272728
$ws.Range("J18").Value = 0.4181818181818182
$ws.Range("K18").Value = 0.09454545454545454
$ws.Range("M18").Value = 0.01818181818181818
$ws.Range("O18").Value = 0.06909090909090909
$ws.Range("S18").Value = 0.1127272727272727
$ws.Range("F19").Value = 0.015625
$ws.Range("H19").Value = 0.18984375
$ws.Range("I19").Value = 0.06796874999999999
$ws.Range("J19").Value = 0.4109375
$ws.Range("K19").Value = 0.10390625
$ws.Range("M19").Value = 0.0203125
$ws.Range("N19").Value = 0.00234375
$ws.Range("O19").Value = 0.065625
$ws.Range("S19").Value = 0.1234375
